$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 190: 四方坪站 data for 2025-12-04 (serial 45995)
$ws.Range("A190").Value = 45995
$ws.Range("B190").Value = "四方坪站充电量(kw)"
$ws.Range("C190").Value = 414.563
$ws.Range("D190").Value = 981.745
$ws.Range("E190").Value = 489.84999999999997
$ws.Range("F190").Value = 264.68
$ws.Range("G190").Value = 247.62699999999998
$ws.Range("H190").Value = 533.858
$ws.Range("I190").Value = 240.33299999999997
$ws.Range("J190").Value = 75.031
$ws.Range("K190").Value = 199.959
$ws.Range("L190").Value = 159.99
$ws.Range("M190").Value = 192.542
$ws.Range("N190").Value = 237.524
$ws.Range("O190").Value = 684.8739999999998
$ws.Range("P190").Value = 1774.3839999999998
$ws.Range("Q190").Value = 604.01
$ws.Range("R190").Value = 451.977
$ws.Range("S190").Value = 397.712
$ws.Range("T190").Value = 232.965
$ws.Range("U190").Value = 74.019
$ws.Range("V190").Value = 93.28
$ws.Range("W190").Value = 105.36
$ws.Range("X190").Value = 104.68
$ws.Range("Y190").Value = 7.3
$ws.Range("Z190").Value = 12.0

# Row 191: 高岭站 data for 2025-12-04 (serial 45995)
$ws.Range("A191").Value = 45995
$ws.Range("B191").Value = "高岭站充电量(kw)"
$ws.Range("C191").Value = 325.5079999999999
$ws.Range("D191").Value = 422.273
$ws.Range("E191").Value = 160.164
$ws.Range("F191").Value = 77.758
$ws.Range("G191").Value = 88.399
$ws.Range("H191").Value = 167.829
$ws.Range("I191").Value = 208.838
$ws.Range("J191").Value = 201.724
$ws.Range("K191").Value = 177.473
$ws.Range("L191").Value = 116.03399999999999
$ws.Range("M191").Value = 204.73899999999998
$ws.Range("N191").Value = 456.0749999999999
$ws.Range("O191").Value = 506.51
$ws.Range("P191").Value = 748.3459999999999
$ws.Range("Q191").Value = 341.884
$ws.Range("R191").Value = 233.593
$ws.Range("S191").Value = 60.507999999999996
$ws.Range("T191").Value = 37.736000000000004
$ws.Range("U191").Value = 141.708
$ws.Range("V191").Value = 42.68
$ws.Range("W191").Value = 33.164
$ws.Range("X191").Value = 61.289
$ws.Range("Y191").Value = 65.531
$ws.Range("Z191").Value = 60.992

# Update view: selection + top-left cell per diff
$ws.Range("D195").Select()
